$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive a plain-numeric-looking string must be forced to Text format
# first (and restored to the default "Normal" style afterward) so Excel does not
# auto-convert them into numeric values, matching the original inlineStr text cells.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D18", "D20", "D22", "D23", "D24", "D28", "D29", "D30", "D31", "D32", "D36", "D37", "D38", "D39", "D40", "D43", "D45", "D46", "D47", "D49")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values exactly as in the new snapshot
$ws.Range('D2').Value = '42.610.05'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '2.289.05'
$ws.Range('E3').Value = '  -3.56%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '301.33'
$ws.Range('E5').Value = '  -2.66%  '
$ws.Range('D6').Value = '98.33'
$ws.Range('E6').Value = '  -6.18%  '
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  -3.30%  '
$ws.Range('D10').Value = '34.54'
$ws.Range('E10').Value = '  -4.08%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = '50.83'
$ws.Range('E12').Value = '  -4.84%  '
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').Value = '6.73'
$ws.Range('E14').Value = '  -3.77%  '
$ws.Range('D15').Value = '2.641.18'
$ws.Range('E15').Value = '  -3.58%  '
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '2.280.67'
$ws.Range('E17').Value = '  -4.00%  '
$ws.Range('D18').Value = '0.793'
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D19').Value = '42.526.21'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').Value = '11.62'
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').Value = '0.0₃0898'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').Value = '6.03'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('D23').Value = '67.11'
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').Value = '235.31'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('E25').Value = '  -4.83%  '
$ws.Range('E26').Value = '  -4.20%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '24.67'
$ws.Range('E28').Value = '  -4.39%  '
$ws.Range('D29').Value = '2.30'
$ws.Range('E29').Value = '  +8.48%  '
$ws.Range('D30').Value = '34.31'
$ws.Range('E30').Value = '  -6.12%  '
$ws.Range('D31').Value = '164.89'
$ws.Range('E31').Value = '  +2.49%  '
$ws.Range('D32').Value = '9.14'
$ws.Range('E32').Value = '  -4.16%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  -4.91%  '
$ws.Range('E35').Value = '  -4.78%  '
$ws.Range('D36').Value = '0.0701'
$ws.Range('E36').Value = '  -5.31%  '
$ws.Range('D37').Value = '4.37'
$ws.Range('E37').Value = '  -6.07%  '
$ws.Range('D38').Value = '2.84'
$ws.Range('E38').Value = '  -8.75%  '
$ws.Range('D39').Value = '16.21'
$ws.Range('E39').Value = '  -11.37%  '
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').Value = '  -7.67%  '
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').Value = '2.40'
$ws.Range('E43').Value = '  -7.04%  '
$ws.Range('D44').Value = '1.970.41'
$ws.Range('E44').Value = '  -3.22%  '
$ws.Range('D45').Value = '0.0284'
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('D46').Value = '18.15'
$ws.Range('E46').Value = '  -7.43%  '
$ws.Range('D47').Value = '9.77'
$ws.Range('E47').Value = '  -7.60%  '
$ws.Range('E48').Value = '  -8.57%  '
$ws.Range('D49').Value = '4.73'
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').Value = '2.513.69'
$ws.Range('E51').Value = '  -3.53%  '

# Restore default styling on the cells we temporarily forced to Text format
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
